$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.932.26"
$ws.Range("E2").Value = "  +3.09%  "

$ws.Range("D3").Value = "2.339.39"
$ws.Range("E3").Value = "  +2.43%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  -1.38%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.85"
$ws.Range("E5").Value = "  +0.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.64"
$ws.Range("E6").Value = "  +6.43%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.634"
$ws.Range("E7").Value = "  +2.00%  "

$ws.Range("E9").Value = "  +4.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.36"
$ws.Range("E10").Value = "  +7.08%  "

$ws.Range("E11").Value = "  +2.29%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.54"
$ws.Range("E12").Value = "  +3.90%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.02"
$ws.Range("E13").Value = "  +4.20%  "

$ws.Range("E14").Value = "  -0.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.51"
$ws.Range("E15").Value = "  +3.15%  "

$ws.Range("D16").Value = "2.690.85"
$ws.Range("E16").Value = "  +2.32%  "

$ws.Range("D17").Value = "2.339.51"
$ws.Range("E17").Value = "  +1.24%  "

$ws.Range("D18").Value = "43.874.17"
$ws.Range("E18").Value = "  +3.60%  "

$ws.Range("E19").Value = "  +3.99%  "

$ws.Range("E20").Value = "  +2.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.96"
$ws.Range("E21").Value = "  -3.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.49"
$ws.Range("E22").Value = "  +2.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.49"
$ws.Range("E23").Value = "  -0.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.88"
$ws.Range("E24").Value = "  +2.24%  "

$ws.Range("E25").Value = "  +5.08%  "

$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.60"
$ws.Range("E27").Value = "  +11.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.16"
$ws.Range("E28").Value = "  +4.85%  "

$ws.Range("E29").Value = "  -0.85%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "39.54"
$ws.Range("E30").Value = "  +9.75%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.59"
$ws.Range("E31").Value = "  +1.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "168.41"
$ws.Range("E32").Value = "  +1.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0904"
$ws.Range("E33").Value = "  +4.84%  "

$ws.Range("E34").Value = "  +9.55%  "

$ws.Range("E35").Value = "  +1.47%  "

$ws.Range("E36").Value = "  +4.47%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.72"
$ws.Range("E37").Value = "  +5.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0364"
$ws.Range("E38").Value = "  +5.24%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.92"
$ws.Range("E39").Value = "  +10.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.79"
$ws.Range("E40").Value = "  +2.81%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.74"
$ws.Range("E41").Value = "  +10.40%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.59"
$ws.Range("E42").Value = "  +10.26%  "

$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.66"
$ws.Range("E43").Value = "  +14.38%  "

$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.240"
$ws.Range("E44").Value = "  +6.45%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "71.86"
$ws.Range("E45").Value = "  +4.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.01"
$ws.Range("E46").Value = "  +0.64%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "115.38"
$ws.Range("E47").Value = "  +4.43%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.659.83"
$ws.Range("E48").Value = "  -2.95%  "

$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.01"
$ws.Range("E49").Value = "  +3.89%  "

$ws.Range("B50").Value = "TheGraph"
$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.217"
$ws.Range("E50").Value = "  +17.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "76.29"
$ws.Range("E51").Value = "  -4.42%  "
